# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets,
# reflecting a refreshed data scrape (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 7045
$ws1.Cells.Item(7, 6).Value = 154
$ws1.Cells.Item(10, 6).Value = 14
$ws1.Cells.Item(11, 6).Value = 53
$ws1.Cells.Item(12, 6).Value = 201
$ws1.Cells.Item(17, 6).Value = 3645
$ws1.Cells.Item(21, 6).Value = 24
$ws1.Cells.Item(23, 6).Value = 2274
$ws1.Cells.Item(25, 6).Value = 257
$ws1.Cells.Item(30, 6).Value = 19
$ws1.Cells.Item(31, 6).Value = 159
$ws1.Cells.Item(32, 6).Value = 282
$ws1.Cells.Item(33, 6).Value = 109

# Sheet "全部类型": row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 7045
$ws4.Cells.Item(8, 6).Value = 154
$ws4.Cells.Item(11, 6).Value = 14
$ws4.Cells.Item(12, 6).Value = 53
$ws4.Cells.Item(13, 6).Value = 201
$ws4.Cells.Item(18, 6).Value = 3645
$ws4.Cells.Item(22, 6).Value = 24
$ws4.Cells.Item(24, 6).Value = 2274
$ws4.Cells.Item(26, 6).Value = 257
$ws4.Cells.Item(31, 6).Value = 19
$ws4.Cells.Item(32, 6).Value = 159
$ws4.Cells.Item(33, 6).Value = 282
$ws4.Cells.Item(34, 6).Value = 109
